# Applies the "Se crea la clase tipo tarea CargarEprepago" change to the
# e-prepago data-driven test workbook:
#   - Datos!B2: numeroDocumento 93221452 -> 93221453
#   - Datos!D2: usuario autotest27 -> autotest28 (new shared string)
#   - Datos sheet view: scroll back to show B2 (drop topLeftCell="I1")
#     and select B2 instead of S3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update the document number used for the new test row.
$ws.Range("B2").Value = 93221453

# Update the "usuario" column to reference the new autotest user.
$ws.Range("D2").Value = "autotest28"

# Bring the sheet view back to the left and select B2 (matches the
# sheetView/selection seen after re-saving with Excel).
$ws.Activate()
$ws.Range("B2").Select()
